$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fix normals problem on cylinder" is fixed now - remove its entire row,
# shifting the rows below it up by one.
$ws.Rows.Item(2).Delete()

# "Vertex welding in model compiler" (now row 3) drops to the bottom of the
# list; everything between it and the end moves up one row to fill the gap.
$taskRow = 3
$lastRow = 6

$movingTask = $ws.Range("A" + $taskRow).Value()
$movingEstimate = $ws.Range("B" + $taskRow).Value()

for ($r = $taskRow; $r -lt $lastRow; $r++) {
    $nextR = $r + 1
    $ws.Range("A" + $r).Value = $ws.Range("A" + $nextR).Value()
    $ws.Range("B" + $r).Value = $ws.Range("B" + $nextR).Value()
}

$ws.Range("A" + $lastRow).Value = $movingTask
$ws.Range("B" + $lastRow).Value = $movingEstimate

# Match the post-edit selection: the whole of row 2 is selected.
$ws.Range("A2:XFD2").Select()
